$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.437.77"
$ws.Range("E2").Value = "  -0.46%  "

$ws.Range("D3").Value = "1.952.33"
$ws.Range("E3").Value = "  +0.46%  "

$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.19"
$ws.Range("E5").Value = "  -0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("E6").Value = "  +2.88%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.91"
$ws.Range("E7").Value = "  +3.25%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.375"
$ws.Range("E9").Value = "  +2.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0781"
$ws.Range("E10").Value = "  -3.06%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  +0.86%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.841"
$ws.Range("E12").Value = "  +3.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.05"
$ws.Range("E13").Value = "  +5.10%  "

$ws.Range("D14").Value = "2.241.55"
$ws.Range("E14").Value = "  +0.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.50"
$ws.Range("E15").Value = "  -2.36%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.24"
$ws.Range("E16").Value = "  +1.15%  "

$ws.Range("D17").Value = "1.964.52"
$ws.Range("E17").Value = "  +1.00%  "

$ws.Range("D18").Value = "36.369.29"
$ws.Range("E18").Value = "  -0.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.01"
$ws.Range("E19").Value = "  -0.48%  "

$ws.Range("D20").Value = "0.0₃0849"
$ws.Range("E20").Value = "  -0.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "228.91"
$ws.Range("E21").Value = "  +0.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.04"
$ws.Range("E22").Value = "  +1.29%  "

$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("E24").Value = "  +1.33%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("E25").Value = "  +1.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.144"
$ws.Range("E26").Value = "  +7.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.07"
$ws.Range("E27").Value = "  -1.45%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.02"
$ws.Range("E28").Value = "  +0.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.17"
$ws.Range("E29").Value = "  -0.39%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.31"
$ws.Range("E30").Value = "  +19.79%  "

$ws.Range("E31").Value = "  +1.65%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.75"
$ws.Range("E32").Value = "  +2.17%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0609"
$ws.Range("E33").Value = "  -1.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.43"
$ws.Range("E34").Value = "  +6.33%  "

$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.25"
$ws.Range("E36").Value = "  +1.76%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.36"
$ws.Range("E37").Value = "  +1.93%  "

$ws.Range("E38").Value = "  -0.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.39"
$ws.Range("E39").Value = "  -12.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0962"
$ws.Range("E40").Value = "  -2.44%  "

$ws.Range("E41").Value = "  -0.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.15"
$ws.Range("E42").Value = "  +0.70%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0208"
$ws.Range("E43").Value = "  -0.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.74"
$ws.Range("E44").Value = "  -0.51%  "

$ws.Range("D45").Value = "1.357.02"
$ws.Range("E45").Value = "  +1.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.13"
$ws.Range("E46").Value = "  +1.85%  "

$ws.Range("E47").Value = "  -0.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.17"
$ws.Range("E48").Value = "  +0.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.83"
$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.74"
$ws.Range("E50").Value = "  +6.22%  "

$ws.Range("D51").Value = "2.136.57"
$ws.Range("E51").Value = "  +0.78%  "

